$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("id", "nome", "email", "evento", "filhos", "idades")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# F1 is a brand-new header cell; give it the same formatting as the rest
# of the header row (bold, centered, thin border) by copying E1's format.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

$data = @(
    @(1, "dsadsada", "laerciogu01@gmail.com", "dsa", "dsadsadas", "dasdsa"),
    @(2, "dsadsadafgdg", "laerciogu01@gmail.com", "dsa", "dsadsadas", "dasdsa"),
    @(3, "fdsfdsfds", "fsdfsfs", "fsdfdsfs", "fsfsdfds", "fsdfds"),
    @(4, "tesnatdo", "dsada", "dasdasdas", "dasdada", "dasdsa"),
    @(5, "dsfdsfsdfsd", "fsdfsdfs", "fsdfsdfsfsd", "fsdfsfsfds", "fsdfs"),
    @(6, "fsdfsdfsd", "fdsfsfs", "fsfdsfs", "fsdfs", "fsdfsfds"),
    @(7, "dsadasda", "dasdsadasdsa", "eeeeeeeeeeeeeee", "eeeeeeeeeeeeeeee", "eeeeeeeeeeeee"),
    @(8, "fsdfsdfsd", "ssssssssssssss", "ssssssssssssssss", "ssssssssssssss", "sssssssssssssss"),
    @(9, "dsdsdsd", "dsdsds", "dsdsds", "dsdsds", "dsdsds")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
